$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.388.16"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.281.48"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.91"
$ws.Range("E5").Value = "  -4.53%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -1.69%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.51"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.13"
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.964"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.06"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "2.626.77"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "2.284.08"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "42.353.83"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.79"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.37"
$ws.Range("E23").Value = "  -6.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.29"
$ws.Range("E24").Value = "  -2.50%  "
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.61"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.05"
$ws.Range("E28").Value = "  +14.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.04"
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.83"
$ws.Range("E31").Value = "  -5.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.26"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0848"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.56"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.48"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.04"
$ws.Range("E42").Value = "  +9.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.09"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").Value = "1.714.95"
$ws.Range("E47").Value = "  +7.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.84"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.76"
$ws.Range("E49").Value = "  -4.55%  "
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.12"
$ws.Range("E51").Value = "  -1.82%  "
